$d = $word.ActiveDocument

# Correction to misclassified observation: "other" -> "unclassified"
# Use whole-word matching to avoid touching "another" elsewhere in the document.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("other", $true, $true, $false, $false, $false, $true, 1, $false, "unclassified", 2)

# Drop the unused built-in Header/Footer (and their linked character) styles.
# Delete in reverse definition order so index-based lookups stay valid.
$d.Styles("FooterChar").Delete()
$d.Styles("Footer").Delete()
$d.Styles("HeaderChar").Delete()
$d.Styles("Header").Delete()
